$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose price text needs to be updated (values are stored as text,
# matching the original inline-string cell type in the sheet).
$targetCells = @("D2","D3","D4","D5","D6","D7","D8","D9","D10","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D23","D24","D25","D40","D41","D44","D45","D47","D48")

$newValues = @{
    "D2" = "247.21"
    "D3" = "22.39"
    "D4" = "5.534"
    "D5" = "0.05636"
    "D6" = "6.469"
    "D7" = "0.8057"
    "D8" = "1.060"
    "D9" = "0.1434"
    "D10" = "0.07309"
    "D12" = "0.1312"
    "D13" = "0.02925"
    "D14" = "0.09263"
    "D15" = "0.001659"
    "D16" = "3.211"
    "D17" = "0.04720"
    "D18" = "0.0005836"
    "D19" = "0.006276"
    "D20" = "0.001053"
    "D21" = "0.004118"
    "D23" = "3.967"
    "D24" = "3.379"
    "D25" = "2.134"
    "D40" = "0.04190"
    "D41" = "0.006877"
    "D44" = "0.009840"
    "D45" = "0.00005633"
    "D47" = "0.6807"
    "D48" = "0.02409"
}

foreach ($addr in $targetCells) {
    $cell = $ws.Range($addr)
    # Force text number-format so the numeric-looking string is not
    # auto-converted to a numeric cell value.
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$addr]
    # Restore default style so no stray formatting/quote-prefix is left behind.
    $cell.Style = "Normal"
}
